$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = 0.8257501125335693
$ws.Range("C2").Value = 0.7060000896453857
$ws.Range("D2").Value = 0.4322500228881836

# Row 3
$ws.Range("B3").Value = 0.8842500448226929
$ws.Range("C3").Value = 0.8457499742507935
$ws.Range("D3").Value = 0.706000030040741

# Row 4 - C4 removed
$ws.Range("B4").Value = 0.8058431148529053
$ws.Range("C4").ClearContents()
$ws.Range("D4").Value = 0.3888919949531555

# Row 5 - C5 removed
$ws.Range("B5").Value = 0.8400057554244995
$ws.Range("C5").ClearContents()
$ws.Range("D5").Value = 0.6442725658416748

# Row 6 - D6 added
$ws.Range("B6").Value = 0.833939790725708
$ws.Range("C6").Value = 0.7052809000015259
$ws.Range("D6").Value = 0.4375340342521667

# Row 7 - D7 added
$ws.Range("B7").Value = 0.8957823514938354
$ws.Range("C7").Value = 0.8444927930831909
$ws.Range("D7").Value = 0.7132723331451416

# Row 8
$ws.Range("B8").Value = 0.8287500143051147
$ws.Range("C8").Value = 0.9390000104904175
$ws.Range("D8").Value = 0.940000057220459

# Row 9
$ws.Range("B9").Value = 0.8282500505447388
$ws.Range("C9").Value = 0.9390000104904175
$ws.Range("D9").Value = 0.8957500457763672

# Row 10 - D10 added
$ws.Range("B10").Value = 0.7939512729644775
$ws.Range("C10").Value = 0.5879067182540894
$ws.Range("D10").Value = 0.5111579895019531

# Row 11 - B11 removed
$ws.Range("B11").ClearContents()
$ws.Range("C11").Value = 0.4135605990886688
$ws.Range("D11").Value = 0.5410357117652893

# Row 12 - D12 added
$ws.Range("B12").Value = 0.7939512729644775
$ws.Range("C12").Value = 0.5879067182540894
$ws.Range("D12").Value = 0.5111579895019531
